# feat(CWL): mod integrity check
# Adds three new rows (60-62) to the localization sheet for the new
# "missing mods" warning dialog (key / JP source / CN translation),
# and tightens the wrap formatting on the now-shorter D59 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) D59 loses its wrap-text formatting (string got shorter) - s=10 -> s=3
#    Grab a same-styled-but-unaffected donor (C59, which already is s=3)
#    and paste its format onto D59, then restore D59's own value.
# ---------------------------------------------------------------------
$ws.Range("C59").Copy()
$ws.Paste($ws.Range("D59"))
$ws.Range("D59").Value = "无法为商人: {0} 加载自定义库存"

# ---------------------------------------------------------------------
# 2) New row 60 - key / JP / CN (CN is rich text: two runs / two fonts)
# ---------------------------------------------------------------------
$ws.Range("A59").Copy()
$ws.Paste($ws.Range("A60"))
$ws.Range("A60").Value = "cwl_warn_missing_mods"

$ws.Range("C59").Copy()
$ws.Paste($ws.Range("C60"))
$ws.Range("C60").Value = "現在のセーブから欠落しているMOD：" + [char]10 + "{0}"

$ws.Range("D58").Copy()
$ws.Paste($ws.Range("D60"))
$d60text = "当前存档中缺失的模组：" + [char]10 + "{0}"
$ws.Range("D60").Value = $d60text
$d60run1 = $ws.Range("D60").Characters(1, 11)
$d60run1.Font.Name = "宋体"
$d60run1.Font.Size = 15.8
$d60run2 = $ws.Range("D60").Characters(12, 4)
$d60run2.Font.Name = "Cascadia Code"
$d60run2.Font.Size = 15.8

$ws.Rows.Item(60).RowHeight = 46.5

# ---------------------------------------------------------------------
# 3) New row 61 - key / JP / CN (plain strings)
# ---------------------------------------------------------------------
$ws.Range("A59").Copy()
$ws.Paste($ws.Range("A61"))
$ws.Range("A61").Value = "cwl_warn_missing_mods_yes"

$ws.Range("C59").Copy()
$ws.Paste($ws.Range("C61"))
$ws.Range("C61").Value = "セーブせずに終了"

$ws.Range("D42").Copy()
$ws.Paste($ws.Range("D61"))
$ws.Range("D61").Value = "不保存并返回至标题"

$ws.Rows.Item(61).RowHeight = 23.25

# ---------------------------------------------------------------------
# 4) New row 62 - key / JP / CN (plain strings)
# ---------------------------------------------------------------------
$ws.Range("A59").Copy()
$ws.Paste($ws.Range("A62"))
$ws.Range("A62").Value = "cwl_warn_missing_mods_no"

$ws.Range("C59").Copy()
$ws.Paste($ws.Range("C62"))
$ws.Range("C62").Value = "プレイを続ける"

$ws.Range("D42").Copy()
$ws.Paste($ws.Range("D62"))
$ws.Range("D62").Value = "继续游玩"

$ws.Rows.Item(62).RowHeight = 23.25

# ---------------------------------------------------------------------
# 5) Selection / view bookkeeping to mirror the saved workbook state
# ---------------------------------------------------------------------
$ws.Range("D65").Select()

$excel.CutCopyMode = $false
